# Generate Report for Archive
#
# The localization status report was regenerated: every cell that used to
# read "Ready for handoff" now reads "In Translation" (Overview!E2:F2 for
# the zh-cn/de-de summary columns, and the per-locale "Status" column on
# the "zh-cn" and "de-de" sheets). The Status-ish columns are narrower now
# that the new text is shorter, so their widths are refreshed to match.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status text wherever it appears ------------------------
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $value = $cell.Value()
        if ($oldStatus -eq $value) {
            $cell.Value = $newStatus
        }
    }
}

# --- Re-fit the columns that held the status text -----------------------
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Columns.Item(5).ColumnWidth = 12.5   # E: zh-cn
$ws1.Columns.Item(6).ColumnWidth = 12.5   # F: de-de

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Columns.Item(3).ColumnWidth = 12.5   # C: Status

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Columns.Item(3).ColumnWidth = 12.5   # C: Status
